$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'263.47"
$ws.Range("D3").Value = "'22.95"
$ws.Range("D4").Value = "'6.185"
$ws.Range("D5").Value = "'0.06245"
$ws.Range("D6").Value = "'6.741"
$ws.Range("D7").Value = "'3.446"
$ws.Range("D8").Value = "'1.345"
$ws.Range("D9").Value = "'0.7959"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01338"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1588"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08124"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03416"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03084"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09333"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.709"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001690"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04796"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").Value = "'0.006220"
$ws.Range("D20").Value = "'0.006186"
$ws.Range("D21").Value = "'0.001094"
$ws.Range("D23").Value = "'3.696"
$ws.Range("D24").Value = "'2.214"
$ws.Range("D25").Value = "'0.3339"
$ws.Range("D27").Value = "'0.0003197"
$ws.Range("D40").Value = "'0.04634"
$ws.Range("D41").Value = "'0.007099"
$ws.Range("D42").Value = "'0.1123"
$ws.Range("D43").Value = "'0.003333"
$ws.Range("D44").Value = "'0.01016"
$ws.Range("D45").Value = "'0.002966"
$ws.Range("D46").Value = "'0.00005885"
$ws.Range("D48").Value = "'0.6990"
$ws.Range("D49").Value = "'0.09097"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
$ws.Range("D50").Value = "'0.00002097"